$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.261.35'
$ws.Range("E2").Value = '  +1.16%  '
$ws.Range("D3").Value = '2.483.32'
$ws.Range("E3").Value = '  +3.17%  '
$ws.Range("E4").Value = '  -0.33%  '
$ws.Range("D5").Value = '577.47'
$ws.Range("E5").Value = '  +0.74%  '
$ws.Range("D6").Value = '146.83'
$ws.Range("E6").Value = '  +0.75%  '
$ws.Range("E7").Value = '  +0.23%  '
$ws.Range("D8").Value = '0.540'
$ws.Range("E8").Value = '  -0.19%  '
$ws.Range("D9").Value = '2.483.08'
$ws.Range("E9").Value = '  +1.98%  '
$ws.Range("D10").Value = '0.111'
$ws.Range("E10").Value = '  +0.58%  '
$ws.Range("E11").Value = '  +1.92%  '
$ws.Range("D12").Value = '5.26'
$ws.Range("E12").Value = '  +0.44%  '
$ws.Range("E13").Value = '  +0.36%  '
$ws.Range("D14").Value = '28.61'
$ws.Range("E14").Value = '  +4.84%  '
$ws.Range("E15").Value = '  +1.49%  '
$ws.Range("D16").Value = '2.935.30'
$ws.Range("E16").Value = '  +3.22%  '
$ws.Range("D17").Value = '63.218.70'
$ws.Range("E17").Value = '  +1.00%  '
$ws.Range("D18").Value = '2.485.02'
$ws.Range("E18").Value = '  +2.45%  '
$ws.Range("D19").Value = '8.24'
$ws.Range("E19").Value = '  +3.90%  '
$ws.Range("E20").Value = '  +0.89%  '
$ws.Range("D21").Value = '329.82'
$ws.Range("E21").Value = '  +0.75%  '
$ws.Range("D22").Value = '2.26'
$ws.Range("E22").Value = '  +10.90%  '
$ws.Range("E23").Value = '  +0.04%  '
$ws.Range("D24").Value = '1.00'
$ws.Range("E24").Value = '  +0.22%  '
$ws.Range("D25").Value = '66.29'
$ws.Range("E25").Value = '  +1.10%  '
$ws.Range("D26").Value = '671.95'
$ws.Range("E26").Value = '  +7.16%  '
$ws.Range("D27").Value = '9.81'
$ws.Range("E27").Value = '  +16.08%  '
$ws.Range("D28").Value = '0.0₂01000'
$ws.Range("E28").Value = '  +1.56%  '
$ws.Range("D29").Value = '2.621.20'
$ws.Range("E29").Value = '  +3.61%  '
$ws.Range("D30").Value = '0.998'
$ws.Range("E30").Value = '  +342.55%  '
$ws.Range("D31").Value = '1.47'
$ws.Range("E31").Value = '  +4.01%  '
$ws.Range("D32").Value = '8.09'
$ws.Range("E32").Value = '  -1.23%  '
$ws.Range("E33").Value = '  +1.37%  '
$ws.Range("E34").Value = '  -3.47%  '
$ws.Range("E35").Value = '  +4.36%  '
$ws.Range("E36").Value = '  +0.33%  '
$ws.Range("D37").Value = '4.80'
$ws.Range("E37").Value = '  +0.84%  '
$ws.Range("E38").Value = '  +1.75%  '
$ws.Range("E39").Value = '  -0.47%  '
$ws.Range("E40").Value = '  +0.86%  '
$ws.Range("D41").Value = '150.79'
$ws.Range("E41").Value = '  -0.24%  '
$ws.Range("D42").Value = '2.71'
$ws.Range("E42").Value = '  -1.84%  '
$ws.Range("D43").Value = '1.76'
$ws.Range("E43").Value = '  +0.17%  '
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("D45").Value = '0.0₆0315'
$ws.Range("E45").Value = '  -28.85%  '
$ws.Range("D46").Value = '156.40'
$ws.Range("E46").Value = '  +8.12%  '
$ws.Range("D47").Value = '15.14'
$ws.Range("E47").Value = '  +7.07%  '
$ws.Range("E48").Value = '  +0.45%  '
$ws.Range("E49").Value = '  +0.13%  '
$ws.Range("E50").Value = '  +1.48%  '
$ws.Range("E51").Value = '  -0.22%  '
